# Auto-generated edits applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.057.54"
$ws.Range("E2").Value = "'  -1.11%  "
$ws.Range("D3").Value = "'1.643.73"
$ws.Range("E3").Value = "'  -1.46%  "
$ws.Range("E4").Value = "'  -0.61%  "
$ws.Range("D5").Value = "'217.57"
$ws.Range("E5").Value = "'  -0.96%  "
$ws.Range("E6").Value = "'  -3.12%  "
$ws.Range("E7").Value = "'  -0.56%  "
$ws.Range("E8").Value = "'  -1.96%  "
$ws.Range("D9").Value = "'0.06276"
$ws.Range("E9").Value = "'  -2.04%  "
$ws.Range("D10").Value = "'20.41"
$ws.Range("E10").Value = "'  -2.38%  "
$ws.Range("D11").Value = "'0.07759"
$ws.Range("E11").Value = "'  -1.20%  "
$ws.Range("D12").Value = "'4.468"
$ws.Range("E12").Value = "'  -2.29%  "
$ws.Range("D13").Value = "'1.640.66"
$ws.Range("E13").Value = "'  -1.83%  "
$ws.Range("D14").Value = "'1.869.62"
$ws.Range("D15").Value = "'0.5576"
$ws.Range("E15").Value = "'  +0.60%  "
$ws.Range("D16").Value = "'0.0₅7984"
$ws.Range("E16").Value = "'  -2.66%  "
$ws.Range("D17").Value = "'64.70"
$ws.Range("E17").Value = "'  -1.80%  "
$ws.Range("D18").Value = "'26.054.95"
$ws.Range("E18").Value = "'  -1.19%  "
$ws.Range("E19").Value = "'  -0.61%  "
$ws.Range("D20").Value = "'4.623"
$ws.Range("E20").Value = "'  -1.55%  "
$ws.Range("D21").Value = "'192.53"
$ws.Range("E21").Value = "'  -0.74%  "
$ws.Range("E22").Value = "'  -2.54%  "
$ws.Range("E23").Value = "'  -1.76%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "'  -0.55%  "
$ws.Range("D25").Value = "'146.60"
$ws.Range("E25").Value = "'  +0.03%  "
$ws.Range("E26").Value = "'  -2.87%  "
$ws.Range("E27").Value = "'  -0.83%  "
$ws.Range("E28").Value = "'  -1.65%  "
$ws.Range("E29").Value = "'  -1.51%  "
$ws.Range("D30").Value = "'0.05633"
$ws.Range("E30").Value = "'  -4.03%  "
$ws.Range("E31").Value = "'  -1.93%  "
$ws.Range("D32").Value = "'3.455"
$ws.Range("E32").Value = "'  -5.27%  "
$ws.Range("D34").Value = "'1.594"
$ws.Range("E34").Value = "'  -0.97%  "
$ws.Range("D35").Value = "'2.788"
$ws.Range("E35").Value = "'  -1.54%  "
$ws.Range("D36").Value = "'2.411"
$ws.Range("E36").Value = "'  -0.41%  "
$ws.Range("D37").Value = "'0.9362"
$ws.Range("E37").Value = "'  -3.77%  "
$ws.Range("D38").Value = "'0.5662"
$ws.Range("E38").Value = "'  -3.19%  "
$ws.Range("D39").Value = "'5.946"
$ws.Range("E39").Value = "'  +1.67%  "
$ws.Range("D40").Value = "'0.01574"
$ws.Range("E40").Value = "'  -1.84%  "
$ws.Range("B41").Value = "'mCoin"
$ws.Range("C41").Value = "'https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D41").Value = "'2.573"
$ws.Range("E41").Value = "'  -0.07%  "
$ws.Range("B42").Value = "'Maker"
$ws.Range("C42").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.050.54"
$ws.Range("E42").Value = "'  -1.47%  "
$ws.Range("B43").Value = "'PaxDollar"
$ws.Range("C43").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.004"
$ws.Range("E43").Value = "'  -0.64%  "
$ws.Range("B44").Value = "'TrustWalletToken"
$ws.Range("C44").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "'0.8407"
$ws.Range("E44").Value = "'  -2.85%  "
$ws.Range("B45").Value = "'Quant"
$ws.Range("C45").Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'102.26"
$ws.Range("E45").Value = "'  -2.46%  "
$ws.Range("B46").Value = "'RocketPoolETH"
$ws.Range("C46").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "'1.781.16"
$ws.Range("E46").Value = "'  -1.45%  "
$ws.Range("B47").Value = "'Aave"
$ws.Range("C47").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'56.86"
$ws.Range("E47").Value = "'  -1.84%  "
$ws.Range("B48").Value = "'Frax"
$ws.Range("C48").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.008"
$ws.Range("E48").Value = "'  -0.55%  "
$ws.Range("B49").Value = "'BabyDogeCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "'0.0₈104"
$ws.Range("E49").Value = "'  -2.01%  "
$ws.Range("B50").Value = "'Cronos"
$ws.Range("C50").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05310"
$ws.Range("E50").Value = "'  +2.76%  "
$ws.Range("B51").Value = "'Mantle"
$ws.Range("C51").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.4324"
$ws.Range("E51").Value = "'  -1.42%  "
